# Updated cryptos list on Tue Oct 29 02:12:47 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns on Sheet1 from the latest
# coinranking.com scrape. Three coins changed rank order since the last run
# (Dai/Litecoin around row 23-24, dogwifhat/PolygonEcosystemToken around row
# 44-45, and ARBITRUM/Optimism/BabyDogeCoin around row 49-51), so their Coin (B)
# and Link (C) columns are refreshed too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.572.04"
$ws.Range("E2").Value = "'  +4.19%  "
$ws.Range("D3").Value = "'2.592.99"
$ws.Range("E3").Value = "'  +4.04%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'604.52"
$ws.Range("E5").Value = "'  +2.86%  "
$ws.Range("D6").Value = "'179.52"
$ws.Range("E6").Value = "'  +2.02%  "
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "'  +1.71%  "
$ws.Range("D9").Value = "'2.588.63"
$ws.Range("E9").Value = "'  +3.83%  "
$ws.Range("E10").Value = "'  +17.92%  "
$ws.Range("E11").Value = "'  +0.37%  "
$ws.Range("D12").Value = "'0.348"
$ws.Range("E12").Value = "'  +3.06%  "
$ws.Range("E13").Value = "'  +1.57%  "
$ws.Range("D14").Value = "'0.0000186"
$ws.Range("E14").Value = "'  +8.74%  "
$ws.Range("D15").Value = "'3.046.44"
$ws.Range("E15").Value = "'  +3.26%  "
$ws.Range("D16").Value = "'26.53"
$ws.Range("E16").Value = "'  +3.37%  "
$ws.Range("D17").Value = "'70.368.77"
$ws.Range("E17").Value = "'  +4.05%  "
$ws.Range("D18").Value = "'2.577.06"
$ws.Range("E18").Value = "'  +2.74%  "
$ws.Range("D19").Value = "'7.81"
$ws.Range("E19").Value = "'  +3.86%  "
$ws.Range("D20").Value = "'11.31"
$ws.Range("E20").Value = "'  +3.53%  "
$ws.Range("D21").Value = "'367.20"
$ws.Range("E21").Value = "'  +4.94%  "
$ws.Range("D22").Value = "'4.19"
$ws.Range("E22").Value = "'  +2.86%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'72.10"
$ws.Range("E23").Value = "'  +1.12%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "'  +0.12%  "
$ws.Range("D25").Value = "'4.40"
$ws.Range("E25").Value = "'  +3.33%  "
$ws.Range("D26").Value = "'1.87"
$ws.Range("E26").Value = "'  +8.44%  "
$ws.Range("D27").Value = "'9.55"
$ws.Range("E27").Value = "'  +5.92%  "
$ws.Range("D28").Value = "'2.721.58"
$ws.Range("E28").Value = "'  +3.75%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("D30").Value = "'0.0₃0946"
$ws.Range("E30").Value = "'  +4.84%  "
$ws.Range("D31").Value = "'525.26"
$ws.Range("E31").Value = "'  +4.67%  "
$ws.Range("D32").Value = "'7.96"
$ws.Range("E32").Value = "'  +2.37%  "
$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "'  +3.38%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "'  +3.69%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "'  -0.18%  "
$ws.Range("D36").Value = "'165.10"
$ws.Range("E36").Value = "'  +0.57%  "
$ws.Range("D37").Value = "'0.121"
$ws.Range("E37").Value = "'  +0.44%  "
$ws.Range("D38").Value = "'19.13"
$ws.Range("E38").Value = "'  +4.46%  "
$ws.Range("D39").Value = "'18.94"
$ws.Range("E39").Value = "'  +1.77%  "
$ws.Range("D40").Value = "'1.37"
$ws.Range("E40").Value = "'  +4.01%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "'  +4.31%  "
$ws.Range("E42").Value = "'  +0.05%  "
$ws.Range("D43").Value = "'4.99"
$ws.Range("E43").Value = "'  +3.18%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.57"
$ws.Range("E44").Value = "'  +3.73%  "
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").Value = "'0.328"
$ws.Range("E45").Value = "'  -0.42%  "
$ws.Range("D46").Value = "'39.17"
$ws.Range("E46").Value = "'  +1.14%  "
$ws.Range("D47").Value = "'153.42"
$ws.Range("E47").Value = "'  +4.36%  "
$ws.Range("D48").Value = "'3.66"
$ws.Range("E48").Value = "'  +2.96%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₆0267"
$ws.Range("E49").Value = "'  +4.24%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'0.531"
$ws.Range("E50").Value = "'  +2.81%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "'1.65"
$ws.Range("E51").Value = "'  +4.75%  "
